$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "About" sheet: update the "last updated" date in C1.
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# ---------------------------------------------------------------------
# 2. Shared-string cleanup: every label of the form
#    "... : NoSettings" on the "BAU Emissions" sheet becomes "... : test".
#    These labels live in column A (rows 1-300).
# ---------------------------------------------------------------------
$wsBau = $wb.Worksheets.Item("BAU Emissions")
$wsBau.Range("A1:A300").Replace(" : NoSettings", " : test")

# ---------------------------------------------------------------------
# 3. "BAU Emissions" sheet: update the data row for
#    "Industrial Sector Energy Related Emissions before CCS[natural gas if,
#    iron and steel 241,CO2]" (row 94, years 2032-2050 in columns M:AE).
# ---------------------------------------------------------------------
$wsBau.Range("M94:AE94").Value = 5005380
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300

# ---------------------------------------------------------------------
# 4. View-state updates.
#    "BAU Emissions" picks up a new selection (A30:AE280, active cell A30).
#    "Current and Planned Capacity" keeps its scroll position but is no
#    longer the selected tab. "About" becomes the selected tab.
# ---------------------------------------------------------------------
$wsBau.Activate()
$wsBau.Range("A30:AE280").Select()

$wsAbout.Activate()

Write-Host "done"
